$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# Row 18 and row 19 swap first/last name, red_shirt, dev_trait, and overall_start
# (the roster order changed for these two WR/physical players).
$ws.Cells.Item(18, 1).Value = "Gene"
$ws.Cells.Item(18, 2).Value = "Kauffman"
$ws.Cells.Item(18, 4).Value = $false
$ws.Cells.Item(18, 10).Value = "elite"
$ws.Cells.Item(18, 11).Value = 80

$ws.Cells.Item(19, 1).Value = "Skylar"
$ws.Cells.Item(19, 2).Value = "Snelling"
$ws.Cells.Item(19, 4).Value = $true
$ws.Cells.Item(19, 10).Value = "star"
$ws.Cells.Item(19, 11).Value = 79

# Fill in overall_end (column L) for every player row.
$overallEnd = @{2=99; 3=86; 4=80; 5=78; 6=78; 7=95; 8=83; 9=82; 10=78; 11=79; 12=75; 13=74; 14=92; 15=86; 16=85; 17=81; 18=80; 19=79; 20=77; 21=75; 22=75; 23=73; 24=83; 25=84; 26=76; 27=74; 28=73; 29=89; 30=80; 31=69; 32=86; 33=83; 34=68; 35=89; 36=81; 37=87; 38=78; 39=76; 40=92; 41=85; 42=80; 43=78; 44=77; 45=88; 46=78; 47=78; 48=72; 49=86; 50=85; 51=75; 52=87; 53=84; 54=80; 55=78; 56=73; 57=84; 58=80; 59=78; 60=70; 61=94; 62=80; 63=79; 64=74; 65=92; 66=89; 67=79; 68=74; 69=73; 70=88; 71=83; 72=83; 73=81; 74=80; 75=78; 76=77; 77=70; 78=87; 79=81; 80=66; 81=87; 82=80; 83=79; 84=78; 85=83; 86=78}
foreach ($row in $overallEnd.Keys) {
    $ws.Cells.Item($row, 12).Value = $overallEnd[$row]
}

# Restore the view: zoomed to 150% with L86 selected.
$excel.ActiveWindow.Zoom = 150
$ws.Range("L86").Select()
